# Updated symbol list on Fri Feb 17 12:00:08 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '309.86'
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '-3.09%'
$c = $ws.Range("G2")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '50.62'
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '3.53%'
$c = $ws.Range("G3")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '5.157'
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '-1.60%'
$c = $ws.Range("G4")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '0.07779'
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '-3.79%'
$c = $ws.Range("G5")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '4.500'
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '-1.96%'
$c = $ws.Range("G6")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '1.346'
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '11.70%'
$c = $ws.Range("G7")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '1.566'
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '-4.52%'
$c = $ws.Range("G8")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '-6.33%'
$c = $ws.Range("G9")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.1982'
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '2.36%'
$c = $ws.Range("G10")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.09622'
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '1.80%'
$c = $ws.Range("G11")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.04744'
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '2.90%'
$c = $ws.Range("G12")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '-0.53%'
$c = $ws.Range("G13")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.001268'
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '-4.62%'
$c = $ws.Range("G14")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.005793'
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '-1.24%'
$c = $ws.Range("G15")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '2,012.63%'
$c = $ws.Range("G16")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '-0.34%'
$c = $ws.Range("G17")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '2.433'
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '0.37%'
$c = $ws.Range("G18")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '1.97%'
$c = $ws.Range("G19")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '8.030'
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '-1.16%'
$c = $ws.Range("G20")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '0.1376'
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '-0.44%'
$c = $ws.Range("G21")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.3095'
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '-0.99%'
$c = $ws.Range("G22")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.04163'
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '-0.24%'
$c = $ws.Range("G23")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '0.001270'
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '-2.69%'
$c = $ws.Range("G24")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.003941'
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '-7.16%'
$c = $ws.Range("G25")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '-0.08%'
$c = $ws.Range("G26")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("G27")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("G28")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("G29")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("G30")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("G31")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("G32")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("G33")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("G34")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("G35")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("G36")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("G37")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.02593'
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '-3.86%'
$c = $ws.Range("G38")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.05997'
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '5.74%'
$c = $ws.Range("G39")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '74.14%'
$c = $ws.Range("G40")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.007884'
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '-0.56%'
$c = $ws.Range("G41")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.1425'
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '-1.17%'
$c = $ws.Range("G42")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.008391'
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '8.81%'
$c = $ws.Range("G43")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.007680'
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '-5.20%'
$c = $ws.Range("G44")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.3388'
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '6.07%'
$c = $ws.Range("G45")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.00007341'
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '6.15%'
$c = $ws.Range("G46")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '-0.16%'
$c = $ws.Range("G47")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '-34.60%'
$c = $ws.Range("G48")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '-20.76%'
$c = $ws.Range("G49")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '-0.16%'
$c = $ws.Range("G50")
$c.NumberFormat = "@"
$c.Value = '12'

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0002000'
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '-0.16%'
$c = $ws.Range("G51")
$c.NumberFormat = "@"
$c.Value = '12'
